$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country list ordering (Chequia/Irlanda swap, Jordania moved)
$ws.Range("A71").Value = "Chequia"
$ws.Range("A72").Value = "Irlanda"
$ws.Range("A136").Value = "Jordania"
$ws.Range("A137").Value = "Bahamas"
$ws.Range("A138").Value = "Aruba"
$ws.Range("A139").Value = "Estonia"

# Update statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes)
$ws.Range("B4").Value = 6523197
$ws.Range("C4").Value = 8966
$ws.Range("D4").Value = 3799897
$ws.Range("E4").Value = 2528886
$ws.Range("G4").Value = 382
$ws.Range("H4").Value = 194414
$ws.Range("B5").Value = 4417550
$ws.Range("C5").Value = 50114
$ws.Range("D5").Value = 3433604
$ws.Range("E5").Value = 909579
$ws.Range("G5").Value = 444
$ws.Range("H5").Value = 74367
$ws.Range("B6").Value = 4179471
$ws.Range("C6").Value = 14347
$ws.Range("E6").Value = 654118
$ws.Range("G6").Value = 602
$ws.Range("H6").Value = 128119
$ws.Range("B14").Value = 427027
$ws.Range("C14").Value = 1486
$ws.Range("D14").Value = 399555
$ws.Range("E14").Value = 15770
$ws.Range("G14").Value = 20
$ws.Range("H14").Value = 11702
$ws.Range("B16").Value = 355219
$ws.Range("C16").Value = 2659
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 41594
$ws.Range("B21").Value = 284943
$ws.Range("C21").Value = 1673
$ws.Range("D21").Value = 254188
$ws.Range("E21").Value = 23918
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 6837
$ws.Range("B22").Value = 281583
$ws.Range("C22").Value = 1434
$ws.Range("D22").Value = 211272
$ws.Range("E22").Value = 34734
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = 35577
$ws.Range("B24").Value = 255267
$ws.Range("C24").Value = 311
$ws.Range("E24").Value = 15258
$ws.Range("B29").Value = 133897
$ws.Range("C29").Value = 149
$ws.Range("D29").Value = 117701
$ws.Range("E29").Value = 7043
$ws.Range("B31").Value = 120846
$ws.Range("C31").Value = 267
$ws.Range("D31").Value = 117746
$ws.Range("E31").Value = 2895
$ws.Range("D53").Value = 56492
$ws.Range("E53").Value = 647
$ws.Range("B65").Value = 41144
$ws.Range("C65").Value = 588
$ws.Range("D65").Value = 28962
$ws.Range("E65").Value = 11086
$ws.Range("G65").Value = 9
$ws.Range("H65").Value = 1096
$ws.Range("B68").Value = 35460
$ws.Range("C68").Value = 104
$ws.Range("D68").Value = 21557
$ws.Range("E68").Value = 13296
$ws.Range("G68").Value = 8
$ws.Range("H68").Value = 607
$ws.Range("B71").Value = 30576
$ws.Range("C71").Value = 699
$ws.Range("D71").Value = 20164
$ws.Range("E71").Value = 9968
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 444
$ws.Range("B72").Value = 30080
$ws.Range("D72").Value = 23364
$ws.Range("E72").Value = 4938
$ws.Range("H72").Value = 1778
$ws.Range("B91").Value = 12080
$ws.Range("C91").Value = 248
$ws.Range("E91").Value = 7983
$ws.Range("G91").Value = 3
$ws.Range("H91").Value = 293
$ws.Range("B92").Value = 11685
$ws.Range("C92").Value = 61
$ws.Range("E92").Value = 2073
$ws.Range("B93").Value = 10704
$ws.Range("C93").Value = 151
$ws.Range("D93").Value = 6284
$ws.Range("E93").Value = 4098
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = 322
$ws.Range("B99").Value = 9108
$ws.Range("C99").Value = 180
$ws.Range("D99").Value = 4640
$ws.Range("E99").Value = 4375
$ws.Range("G99").Value = 2
$ws.Range("H99").Value = 93
$ws.Range("B100").Value = 8899
$ws.Range("C100").Value = 39
$ws.Range("D100").Value = 7683
$ws.Range("E100").Value = 1145
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 71
$ws.Range("E101").Value = 2553
$ws.Range("G101").Value = 2
$ws.Range("H101").Value = 31
$ws.Range("D103").Value = 7500
$ws.Range("E103").Value = 593
$ws.Range("B133").Value = 3142
$ws.Range("C133").Value = 2
$ws.Range("E133").Value = 184
$ws.Range("B136").Value = 2659
$ws.Range("C136").Value = 78
$ws.Range("D136").Value = 1919
$ws.Range("E136").Value = 721
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 19
$ws.Range("B137").Value = 2657
$ws.Range("C137").Value = 72
$ws.Range("D137").Value = 1088
$ws.Range("E137").Value = 1506
$ws.Range("G137").Value = 4
$ws.Range("H137").Value = 63
$ws.Range("B138").Value = 2589
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 1293
$ws.Range("E138").Value = 1281
$ws.Range("H138").Value = 15
$ws.Range("B139").Value = 2585
$ws.Range("C139").Value = 22
$ws.Range("D139").Value = 2213
$ws.Range("E139").Value = 308
$ws.Range("H139").Value = 64
$ws.Range("B141").Value = 2493
$ws.Range("C141").Value = 102
$ws.Range("D141").Value = 749
$ws.Range("E141").Value = 1705
$ws.Range("B142").Value = 2416
$ws.Range("C142").Value = 70
$ws.Range("E142").Value = 1090
$ws.Range("B144").Value = 2242
$ws.Range("C144").Value = 29
$ws.Range("E144").Value = 409
$ws.Range("B150").Value = 1889
$ws.Range("C150").Value = 180
$ws.Range("D150").Value = 553
$ws.Range("E150").Value = 1324
$ws.Range("B161").Value = 1313
$ws.Range("C161").Value = 2
$ws.Range("E161").Value = 37

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 9 de Septiembre de 2020 a las 18:19"
